# Regenerate merged AHB files
# - For the "group header" rows (first row of each new lfd. Position group,
#   columns B value changes), switch the plain border-only style (s=5) to
#   the grey-filled header style already used by earlier groups: s=2 for
#   most columns, s=3 (bold) for column B. We reuse the already-correct
#   formatting from row 2 (A2:V2) as the template and paste only formats.
# - For every data row in the block (82-173, excluding the untouched
#   "NEU" marker row 167), the "Änderung" column L loses its "ÄNDERUNG"
#   marker text and becomes a plain empty grey-centered cell (style s=4),
#   matching the template already used by L2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$headerRows = @(82, 85, 89, 95, 100, 104, 108, 117, 123, 145, 168, 171)

$templateRow = $ws.Range("A2:V2")

foreach ($r in $headerRows) {
    $templateRow.Copy()
    $dst = $ws.Range("A" + $r + ":V" + $r)
    $dst.PasteSpecial($xlPasteFormats)
}

$templateL = $ws.Range("L2")

$firstRow = 82
$lastRow = 173
$skipRow = 167

for ($r = $firstRow; $r -le $lastRow; $r++) {
    if ($r -eq $skipRow) {
        continue
    }
    $templateL.Copy()
    $cell = $ws.Range("L" + $r)
    $cell.PasteSpecial($xlPasteFormats)
    $cell.ClearContents()
}
